$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G14").Value = "Pass"
$ws.Range("G15").Value = "Pass"
$ws.Range("F14").Value = "It should display the New list of :                                                    *Non-food invoice                                                                                  *Image not clear/cut off - please retake                                                       *Duplicate image upload - no action required                                 *Wrong outlet/section                                                                       *Missing page                                                                                      *Others (Please specify) and Type up to 40 characters"
